$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Property1" (sheet1.xml): insert a new row above the old row 8
# ("Upload") labelled "Force", pushing "Upload" to row 9 and the header row
# ("Desc") to row 10. The new row copies the formatting of the row below it
# (the "Upload" row) and is filled with FALSE booleans, like every other
# boolean-flag row on this sheet.
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Property1")
$ws1.Activate()

$ws1.Range("A8:AT8").Insert(-4121)          # xlShiftDown
$ws1.Range("A9:AT9").Copy()
$ws1.Range("A8:AT8").PasteSpecial(-4122)    # xlPasteFormats
$excel.CutCopyMode = $false

$ws1.Range("A8").Value2 = "Force"
$ws1.Range("B8:AT8").Value2 = $false

# Re-home the frozen pane one row further down (it used to freeze the first
# 9 rows; now it freezes the first 10) and restore the scroll position.
$win1 = $excel.ActiveWindow
$win1.FreezePanes = $false
$ws1.Range("A11").Select()
$win1.FreezePanes = $true
$ws1.Range("A9").Select()

# ---------------------------------------------------------------------------
# Sheet "Property2" (sheet2.xml): same row-insertion pattern, no frozen pane.
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Property2")
$ws2.Activate()

$ws2.Range("A8:AD8").Insert(-4121)          # xlShiftDown
$ws2.Range("A9:AD9").Copy()
$ws2.Range("A8:AD8").PasteSpecial(-4122)    # xlPasteFormats
$excel.CutCopyMode = $false

$ws2.Range("A8").Value2 = "Force"
$ws2.Range("B8:AD8").Value2 = $false

$ws2.Range("A9").Select()

# ---------------------------------------------------------------------------
# Sheet "Record_CommValue" (sheet3.xml): selection moves from B5 to the
# whole of row 7.
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("Record_CommValue")
$ws3.Activate()
$ws3.Rows.Item(7).Select()

# Keep "Record_CommValue" as the active/selected tab, matching the original
# workbook's activeTab.
$ws3.Activate()
